# feat(stats): Actual statistical analysis
#
# Rework the ad-hoc "t-test vs fictional no-change group" scratch columns
# (K:R) on Sheet1 into a real year-over-year statistics table on a new
# Sheet2 (Year / Mean NDVI / Difference / Relative Difference), and strip
# the old scratch formulas/labels off Sheet1 (keeping just the two
# "Additive Change" / "Multiplicative Change" ratio columns that feed the
# chart).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Sheet1 cleanup: remove the old N/STDEV/MEAN/.../CONFIDENCE scratch
#    table (columns K:R) in both blocks, plus the stray commentary cells,
#    keeping the "Additive Change"/"Multiplicative Change" label cells and
#    the J-column change/ratio formulas untouched.
# ---------------------------------------------------------------------
$ws1.Range("K1:R1").ClearContents() | Out-Null
$ws1.Range("K2:R2").ClearContents() | Out-Null
$ws1.Range("N3").ClearContents() | Out-Null

$ws1.Range("K23:R23").ClearContents() | Out-Null
$ws1.Range("K24:R24").ClearContents() | Out-Null
$ws1.Range("N25").ClearContents() | Out-Null

# ---------------------------------------------------------------------
# 2. New Sheet2: per-year mean NDVI, year-over-year difference and
#    relative (%) difference, plus a couple of summary cells.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)

$ws2.Cells.Item(1,1).Value = "Year"
$ws2.Cells.Item(1,2).Value = "Mean NDVI"
$ws2.Cells.Item(1,3).Value = "Difference"
$ws2.Cells.Item(1,4).Value = "Relative Difference"

$years = 2000,2001,2002,2003,2004,2005,2006,2007,2008,2009,2010,2011,2012,2013,2014,2015,2016,2017,2018,2019
$means = 2745.8646153846148,3054.455384615384,3558.9092307692308,3723.956923076923,2871.5907692307692,3228.8830769230772,3156.752307692308,3133.476923076923,2819.1323076923081,2559.4907692307688,2745.8646153846148,2454.5061538461541,2678.2415384615379,3092.104615384615,2581.4907692307688,3574.941538461539,2973.3876923076919,2865.6569230769228,2557.666153846154,3092.2184615384622

for ($i = 0; $i -lt $years.Length; $i++) {
    $row = $i + 2
    $ws2.Cells.Item($row, 1).Value = $years[$i]
    $ws2.Cells.Item($row, 2).Value = $means[$i]
    $ws2.Cells.Item($row, 2).NumberFormat = "0"
}

# Year-over-year absolute difference (shared formula over C4:C21, C3 set
# individually since it is the first data row).
$ws2.Range("C3").Formula = "=B3-B2"
$ws2.Range("C3").NumberFormat = "0"
$ws2.Range("C4:C21").Formula = "=B4-B3"
$ws2.Range("C4:C21").NumberFormat = "0"

# Year-over-year relative difference (percent).
$ws2.Range("D3").Formula = "=C3/B2"
$ws2.Range("D3").NumberFormat = "0%"
$ws2.Range("D4:D21").Formula = "=C4/B3"
$ws2.Range("D4:D21").NumberFormat = "0%"

# Summary: average & stdev of the year-over-year differences, and their
# ratio (coefficient-of-variation-style sanity check).
$ws2.Range("F3").Formula = "=AVERAGE(C3:C21)"
$ws2.Range("F3").NumberFormat = "0.00"
$ws2.Range("F4").Formula = "=_xlfn.STDEV.S(C3:C21)"
$ws2.Range("G4").Formula = "=F3/F4"

$ws2.Columns.Item(2).ColumnWidth = 11.83

# ---------------------------------------------------------------------
# 3. Selections / active tab: Sheet2 becomes the active sheet (selected
#    cell G4), Sheet1 keeps a selection at P18 for when the user flips
#    back to it.
# ---------------------------------------------------------------------
$ws1.Range("P18").Select() | Out-Null
$ws2.Range("G4").Select() | Out-Null
